$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1264.704
$ws.Range("D2").Value = 328.806
$ws.Range("C3").Value = 3355.414
$ws.Range("D3").Value = 442.6559999999999
$ws.Range("C4").Value = 2284.951
$ws.Range("D4").Value = 332.913
$ws.Range("C7").Value = 6778.589999999999
$ws.Range("C11").Value = 1207.721
$ws.Range("C12").Value = 1854.017
$ws.Range("C13").Value = 1491.668
$ws.Range("C14").Value = 2194.491
$ws.Range("C16").Value = 3012.214999999999
$ws.Range("C23").Value = 1539.889
$ws.Range("C24").Value = 2701.49
$ws.Range("C25").Value = 2064.289
$ws.Range("C29").Value = 968.359
$ws.Range("C30").Value = 1315.37
$ws.Range("C31").Value = 1033.525
$ws.Range("C35").Value = 702.9200000000001
$ws.Range("C36").Value = 1225.508
$ws.Range("C37").Value = 932.037
$ws.Range("C48").Value = 660.373
$ws.Range("C51").Value = 604.2430000000001
$ws.Range("C52").Value = 1420.097
$ws.Range("C53").Value = 1461.59
$ws.Range("C54").Value = 424.0599999999999
$ws.Range("D54").Value = 201.897
$ws.Range("C55").Value = 1003.916
$ws.Range("D55").Value = 458.6189999999999
$ws.Range("C56").Value = 598.02
$ws.Range("C63").Value = 1179.992
$ws.Range("C64").Value = 1843.662
$ws.Range("C65").Value = 1517.764
$ws.Range("C66").Value = 1399.235
$ws.Range("C67").Value = 2473.991
$ws.Range("C68").Value = 1864.623
$ws.Range("C69").Value = 885.019
$ws.Range("D69").Value = 215.836
$ws.Range("C70").Value = 1717.05
$ws.Range("D70").Value = 309.674
$ws.Range("C71").Value = 1229.62
$ws.Range("D71").Value = 248.954
$ws.Range("C72").Value = 181.624
$ws.Range("D72").Value = 152.386
$ws.Range("C73").Value = 327.471
$ws.Range("D73").Value = 169.427
$ws.Range("C74").Value = 223.623
$ws.Range("D74").Value = 149.788
$ws.Range("C81").Value = 1729.038
$ws.Range("C82").Value = 3476.469999999999
$ws.Range("C83").Value = 2575.789000000001
$ws.Range("C87").Value = 1360.578
$ws.Range("C88").Value = 2245.711
$ws.Range("C89").Value = 1731.692
$ws.Range("C90").Value = 1477.29
$ws.Range("C91").Value = 2512.252
$ws.Range("C92").Value = 1886.03
$ws.Range("C93").Value = 1341.344
$ws.Range("D93").Value = 205.355
$ws.Range("C94").Value = 2072.827
$ws.Range("C95").Value = 1693.84
$ws.Range("D95").Value = 237.504
